$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-17 11:04:49"
$wsZhCn.Range("H3").Value = "2016-03-17 11:05:19"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-17 11:04:55"
$wsDeDe.Range("H3").Value = "2016-03-17 11:05:24"
